# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values to reflect the repulled data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -9
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -8
$ws.Range("F7").Value = -11
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = -2
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = 0
